$wb = $excel.ActiveWorkbook

# Sheet "具有相當價值之財產" (property of considerable value) is the 6th sheet.
$ws = $wb.Worksheets.Item(6)
$stock = $wb.Worksheets.Item(5)

# --- Row 1 (header labels) ---
# Extend the header row from B1:E1 out to L1, reusing the bold/bordered
# header style already applied to the existing header cells.
$ws.Range("D1").Copy($ws.Range("F1:L1"))

$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "quantity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "property_category"
$ws.Range("G1").Value = "category"
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"
$ws.Range("K1").Value = "source_file"
$ws.Range("L1").Value = "index"

# --- Row 2 (101 - 手錶) ---
$ws.Range("F2").Value = "手錶珠寶"
$ws.Range("G2").Value = "normal"
# "2013-12-26" looks like a date, so a plain .Value write would get parsed
# and stored as a date serial instead of literal text. Pull the same text
# already stored as a shared string on the 股票 (stock) sheet's date
# column and paste only the *value* across, which preserves its text type
# without touching number formats / introducing new cell styles.
$stock.Range("J2").Copy()
$ws.Range("H2").PasteSpecial(-4163)
$ws.Range("I2").Value = "丁守中"
$ws.Range("J2").Value = 515
$ws.Range("K2").Value = "tmpc7fb1"
$ws.Range("L2").Value = 101

# --- Row 3 (102 - was 手錶珠寶, relabelled otherbonds) ---
$ws.Range("B3").Value = "otherbonds"
$ws.Range("F3").Value = "手錶珠寶"
$ws.Range("G3").Value = "normal"
$stock.Range("J2").Copy()
$ws.Range("H3").PasteSpecial(-4163)
$ws.Range("I3").Value = "丁守中"
$ws.Range("J3").Value = 515
$ws.Range("K3").Value = "tmpc7fb1"
$ws.Range("L3").Value = 102
